$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (theta_se) standard-error values, columns B:L
$row4 = @("(0.42)", "(0.41)", "(0.49)", "(1.0)", "(0.44)", "(0.1)", "(0.24)", "(0.76)", "(0.04)", "(0.92)", "(0.38)")

# Row 6 (lambda_se) standard-error values, columns B:L
$row6 = @("(0.01)", "(0.26)", "(0.15)", "(0.3)", "(0.02)", "(0.35)", "(0.54)", "(0.08)", "(0.37)", "(0.6)", "(0.54)")

for ($i = 0; $i -lt $row4.Length; $i++) {
    $col = 2 + $i   # column B = 2
    $ws.Cells.Item(4, $col).Value = $row4[$i]
    $ws.Cells.Item(6, $col).Value = $row6[$i]
}
